# "Generate Report for Handoff"
# Update status from "In Translation" to "Ready for handoff" and refresh the
# associated timestamps on all three report sheets, then widen the status
# columns so the new, longer status text fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-30 14:51:10"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-30 14:50:57"

# --- de-de sheet ------------------------------------------------------------
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-30 14:51:10"

# --- Widen the status columns to fit the new "Ready for handoff" text ------
$overview.Range("E1").EntireColumn.ColumnWidth = 16.33
$overview.Range("F1").EntireColumn.ColumnWidth = 16.33
$zhcn.Range("C1").EntireColumn.ColumnWidth = 16.33
$dede.Range("C1").EntireColumn.ColumnWidth = 16.33
